$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 15: G15 = "25/10/2018", H15 = "31/10/2018"
$ws.Range("G15").Value = "25/10/2018"
$ws.Range("H15").Value = "31/10/2018"

# Update selection to H15
$ws.Range("H15").Select()
